$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for years 2000-2009 (original rows 2-11). This shifts
# the 2010-2020 data (originally rows 12-22) up to rows 2-12.
$ws.Range("A2:E11").EntireRow.Delete()

# Append the new 2021 row (row 13) with only 孤儿数 (D) and
# 家庭儿童收养登记总数 (E) populated (中国公民收养登记数/B and
# 外国公民收养登记数/C have no data yet for 2021).
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 4).Value = 172716
$ws.Cells.Item(13, 5).Value = 12447

# Match the formatting used by the rest of column A (bold, bordered,
# centered header-style cell) by copying it from the row above.
$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4122)
